$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to the new "custom accuracy" (2 decimal place) readings
$ws.Range("B5").Value = 4.8
$ws.Range("C5").Value = 3.27
$ws.Range("D5").Value = 0.72
$ws.Range("E5").Value = 10.14
$ws.Range("F5").Value = 8.39
$ws.Range("G5").Value = 3.78
$ws.Range("H5").Value = 19.2
$ws.Range("I5").Value = 5.82
$ws.Range("J5").Value = 2.5
$ws.Range("K5").Value = 3.7
$ws.Range("L5").Value = 4.16
$ws.Range("M5").Value = 4.23
$ws.Range("N5").Value = 1.22
$ws.Range("O5").Value = 3.76
$ws.Range("P5").Value = 5.31
$ws.Range("Q5").Value = 3.35
$ws.Range("R5").Value = 0.73
$ws.Range("S5").Value = 0.37
$ws.Range("T5").Value = 49.94
$ws.Range("U5").Value = 10.74
$ws.Range("V5").Value = 3.47
$ws.Range("W5").Value = 7.05
$ws.Range("X5").Value = 3.84
$ws.Range("Y5").Value = 0.39
$ws.Range("Z5").Value = 8.69
$ws.Range("AA5").Value = 3.07
$ws.Range("AB5").Value = 2.84
$ws.Range("AC5").Value = 3.31
$ws.Range("AD5").Value = 4.3
$ws.Range("AE5").Value = 0.56
$ws.Range("AF5").Value = 17.53
$ws.Range("AG5").Value = 1.87
$ws.Range("AH5").Value = 4.34

# Remove row 6 entirely (shrinks the used range to A1:AH5)
$ws.Rows.Item(6).Delete()
